$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1117.4546
$ws.Range("I12").Value = 1129.2
$ws.Range("K12").Value = 1129.2
$ws.Range("M12").Value = -959.2
$ws.Range("H17").Value = 2913.5833
$ws.Range("I17").Value = 2899.5
$ws.Range("J17").Value = 2984
$ws.Range("K17").Value = 8698.5
$ws.Range("L17").Value = 8952
$ws.Range("M17").Value = -8530.5
$ws.Range("N17").Value = -9288
$ws.Range("H18").Value = 999
$ws.Range("I18").Value = 999
$ws.Range("K18").Value = 999
$ws.Range("M18").Value = -715
$ws.Range("H33").Value = 851.3
$ws.Range("J33").Value = 997
$ws.Range("L33").Value = 997
$ws.Range("N33").Value = -1455
$ws.Range("H75").Value = 45000
$ws.Range("J75").Value = 45000
$ws.Range("L75").Value = 45000
$ws.Range("N75").Value = -46872
$ws.Range("H78").Value = 45000
$ws.Range("J78").Value = 45000
$ws.Range("L78").Value = 135000
$ws.Range("N78").Value = -144360
$ws.Range("H86").Value = 1199.5
$ws.Range("J86").Value = 1199.5
$ws.Range("L86").Value = 1199.5
$ws.Range("N86").Value = -3445.5
$ws.Range("H89").Value = 1199.5
$ws.Range("J89").Value = 1199.5
$ws.Range("L89").Value = 5997.5
$ws.Range("N89").Value = -17229.5
$ws.Range("H100").Value = 3448.5
$ws.Range("I100").Value = 3699.5
$ws.Range("J100").Value = 3323
$ws.Range("K100").Value = 3699.5
$ws.Range("L100").Value = 3323
$ws.Range("M100").Value = -3158.5
$ws.Range("N100").Value = -4405
$ws.Range("H107").Value = 623.5217
$ws.Range("J107").Value = 497.25
$ws.Range("L107").Value = 497.25
$ws.Range("N107").Value = -4337.25
$ws.Range("H112").Value = 1766.6666
$ws.Range("J112").Value = 2000
$ws.Range("L112").Value = 6000
$ws.Range("N112").Value = -8216
$ws.Range("H113").Value = 4999
$ws.Range("I113").Value = 4999
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4999
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = -1745

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 3105
$ws.Range("I29").Value = 3200
$ws.Range("J29").Value = 3010
$ws.Range("K29").Value = 3200
$ws.Range("L29").Value = 3010
$ws.Range("M29").Value = -2892
$ws.Range("N29").Value = -3626
$ws.Range("H45").Value = 1871.7
$ws.Range("I45").Value = 1886.6666
$ws.Range("K45").Value = 1886.6666
$ws.Range("M45").Value = -1509.6666
$ws.Range("H97").Value = 540.8333
$ws.Range("I97").Value = 509
$ws.Range("K97").Value = 509
$ws.Range("M97").Value = -13
$ws.Range("H102").Value = 1461.4706
$ws.Range("I102").Value = 1223
$ws.Range("J102").Value = 3250
$ws.Range("K102").Value = 1223
$ws.Range("L102").Value = 3250
$ws.Range("M102").Value = 399
$ws.Range("N102").Value = -6494
$ws.Range("H122").Value = 1837.25
$ws.Range("I122").Value = 1837.25
$ws.Range("K122").Value = 5511.75
$ws.Range("M122").Value = -3061.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1003.7
$ws.Range("I20").Value = 805.8570999999999
$ws.Range("K20").Value = 805.8570999999999
$ws.Range("M20").Value = -558.8570999999999
$ws.Range("H99").Value = 3074.8215
$ws.Range("I99").Value = 3181.3333
$ws.Range("K99").Value = 3181.3333
$ws.Range("M99").Value = -1683.3333
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("N124").Value = 0
$ws.Range("H125").Value = 150000
$ws.Range("J125").Value = 150000
$ws.Range("L125").Value = 150000
$ws.Range("N125").Value = -159840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 806.06665
$ws.Range("I22").Value = 830.0769
$ws.Range("J22").Value = 650
$ws.Range("K22").Value = 830.0769
$ws.Range("L22").Value = 650
$ws.Range("M22").Value = -480.0769
$ws.Range("N22").Value = -1350
$ws.Range("H31").Value = 2639.6667
$ws.Range("I31").Value = 1903.9375
$ws.Range("J31").Value = 4111.125
$ws.Range("K31").Value = 1903.9375
$ws.Range("L31").Value = 4111.125
$ws.Range("M31").Value = -1608.9375
$ws.Range("N31").Value = -4701.125
$ws.Range("H34").Value = 2639.6667
$ws.Range("I34").Value = 1903.9375
$ws.Range("J34").Value = 4111.125
$ws.Range("K34").Value = 1903.9375
$ws.Range("L34").Value = 4111.125
$ws.Range("M34").Value = -1701.9375
$ws.Range("N34").Value = -4515.125
$ws.Range("H86").Value = 6530.6
$ws.Range("I86").Value = 3788.25
$ws.Range("K86").Value = 3788.25
$ws.Range("M86").Value = -2665.25
$ws.Range("H89").Value = 6530.6
$ws.Range("I89").Value = 3788.25
$ws.Range("K89").Value = 18941.25
$ws.Range("M89").Value = -13325.25
$ws.Range("H122").Value = 4480.1816
$ws.Range("I122").Value = 1438.8334
$ws.Range("J122").Value = 8129.8
$ws.Range("K122").Value = 4316.5002
$ws.Range("L122").Value = 24389.4
$ws.Range("M122").Value = -1866.5002
$ws.Range("N122").Value = -29289.4
$ws.Range("H134").Value = 1865
$ws.Range("I134").Value = 2000
$ws.Range("K134").Value = 6000
$ws.Range("M134").Value = -3465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 20.25
$ws.Range("I7").Value = 20
$ws.Range("J7").Value = 21
$ws.Range("K7").Value = 60
$ws.Range("L7").Value = 63
$ws.Range("M7").Value = 52
$ws.Range("N7").Value = -287
$ws.Range("H23").Value = 250
$ws.Range("J23").Value = 250
$ws.Range("L23").Value = 750
$ws.Range("N23").Value = -1220
$ws.Range("H50").Value = 626.6667
$ws.Range("I50").Value = 626.6667
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 1880.0001
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = -1399.0001
$ws.Range("H53").Value = 626.6667
$ws.Range("I53").Value = 626.6667
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 1880.0001
$ws.Range("L53").Value = 0
$ws.Range("N53").Value = -1399.0001
$ws.Range("H64").Value = 416.66666
$ws.Range("I64").Value = 416.66666
$ws.Range("K64").Value = 1249.99998
$ws.Range("M64").Value = -979.9999800000001
$ws.Range("H67").Value = 416.66666
$ws.Range("I67").Value = 416.66666
$ws.Range("K67").Value = 1249.99998
$ws.Range("M67").Value = -313.9999800000001
$ws.Range("H113").Value = 448.77777
$ws.Range("J113").Value = 196.33333
$ws.Range("L113").Value = 588.99999
$ws.Range("N113").Value = -4928.99999
$ws.Range("H121").Value = 615.0769
$ws.Range("I121").Value = 168.33333
$ws.Range("J121").Value = 998
$ws.Range("K121").Value = 504.99999
$ws.Range("L121").Value = 2994
$ws.Range("M121").Value = 805.00001
$ws.Range("N121").Value = -5614
$ws.Range("H129").Value = 1256.1666
$ws.Range("I129").Value = 679.6667
$ws.Range("J129").Value = 1832.6666
$ws.Range("K129").Value = 2039.0001
$ws.Range("L129").Value = 5497.9998
$ws.Range("M129").Value = 2960.9999
$ws.Range("N129").Value = -15497.9998
$ws.Range("H131").Value = 3387.4
$ws.Range("I131").Value = 975
$ws.Range("K131").Value = 2925
$ws.Range("M131").Value = 2115

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6382.231
$ws.Range("J122").Value = 6274.778
$ws.Range("L122").Value = 18824.334
$ws.Range("N122").Value = -23724.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2825
$ws.Range("I100").Value = 2825
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2825
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = -2284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5161.3335
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 5161.3335
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("H69").Value = 20044.166
$ws.Range("J69").Value = 20044.166
$ws.Range("L69").Value = 20044.166
$ws.Range("N69").Value = -21542.166
$ws.Range("H72").Value = 20044.166
$ws.Range("J72").Value = 20044.166
$ws.Range("L72").Value = 60132.49800000001
$ws.Range("N72").Value = -67620.49800000001
$ws.Range("H100").Value = 391
$ws.Range("I100").Value = 391
$ws.Range("K100").Value = 782
$ws.Range("M100").Value = -241
$ws.Range("H103").Value = 50000
$ws.Range("J103").Value = 50000
$ws.Range("L103").Value = 50000
$ws.Range("N103").Value = -52344
$ws.Range("H125").Value = 30715
$ws.Range("J125").Value = 30715
$ws.Range("L125").Value = 30715
$ws.Range("N125").Value = -40555
$ws.Range("H132").Value = 1757
$ws.Range("I132").Value = 1575
$ws.Range("J132").Value = 2485
$ws.Range("K132").Value = 4725
$ws.Range("L132").Value = 7455
$ws.Range("M132").Value = -2195
$ws.Range("N132").Value = -12515
